# Apply updated cryptocurrency price/volume data per the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "61.474.38"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -3.53%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.478.73"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -6.06%  "

$ws.Range("E4").Value = "  +0.06%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "553.00"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -4.65%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "146.59"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -5.65%  "

$ws.Range("E7").Value = "  +0.04%  "

$ws.Range("E8").Value = "  -3.25%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.477.66"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -6.00%  "

$ws.Range("E10").Value = "  -8.36%  "

$ws.Range("E11").Value = "  -5.78%  "

$ws.Range("E12").Value = "  -1.45%  "

$ws.Range("E13").Value = "  -6.38%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.19"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -7.62%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.926.37"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -6.03%  "

$ws.Range("E16").Value = "  -8.76%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "61.382.15"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -3.58%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.479.68"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -6.09%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.18"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -7.71%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.06"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -7.77%  "

$ws.Range("E21").Value = "  -7.10%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "321.84"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -6.49%  "

$ws.Range("E23").Value = "  -0.01%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.87"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.45%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "64.09"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -5.76%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0₃0997"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -8.63%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.607.58"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -5.75%  "

$ws.Range("E28").Value = "  -5.49%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "544.12"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -9.72%  "

$ws.Range("E30").Value = "  +0.07%  "

$ws.Range("E31").Value = "  -9.28%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.81"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.55%  "

$ws.Range("E33").Value = "  -5.28%  "

$ws.Range("E34").Value = "  -7.36%  "

$ws.Range("E35").Value = "  -8.07%  "

$ws.Range("E36").Value = "  -10.04%  "

$ws.Range("E37").Value = "  -9.94%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.00"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.06%  "

$ws.Range("E39").Value = "  -4.86%  "

$ws.Range("E40").Value = "  -5.47%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "146.46"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.97%  "

$ws.Range("E42").Value = "  -8.09%  "

$ws.Range("E43").Value = "  +0.04%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "40.44"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.50%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.38"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -6.31%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "147.59"
$ws.Range("D46").Style = "Normal"

$ws.Range("E47").Value = "  -6.49%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "21.31"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -11.79%  "

$ws.Range("E49").Value = "  -7.45%  "

$ws.Range("E50").Value = "  -5.41%  "

$ws.Range("E51").Value = "  -4.94%  "
